# "try to fix render markdown"
#
# The GCS coverage-level description strings (column E, rows 26-73) used a
# literal "<br>" as a line separator, which doesn't render as markdown.
# Replace each "<br>" with a real newline so the text wraps as intended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E26:E37").Value = "Supports the GCS if coverage is **Low** `nOther members: Global South + EU `n(25-33% of world emissions)"
$ws.Range("E38:E49").Value = "Supports the GCS if coverage is **Mid** `nGlobal South + China `n(56% of world emissions)"
$ws.Range("E50:E61").Value = "Supports the GCS if coverage is **High** `nGlobal South + China + EU + various HICs `n(UK, Japan, Korea, Canada...; 64-72% of emissions)"
$ws.Range("E62:E73").Value = "Supports the GCS if coverage is **High**, **color** variant `nGlobal South + China + EU + various HICs `n+ Distributive effects shown using colors on world map"

# Re-fit the affected rows so the newly multi-line text doesn't leave a
# stray custom row height behind (keeps the rest of the sheet untouched).
$ws.Rows("26:73").AutoFit()
